$d = $word.ActiveDocument

# Update the date heading paragraph (not part of the table).
$d.Content.Find.Execute("2026-01-08 Thursday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2026-01-09 Friday", 2)

# The practice table holds the multiplication problems. Cells are
# addressed directly (row/col, 1-based) rather than via document-wide
# Find/Replace because several answers repeat verbatim (e.g. "456x6=2736"
# appears twice with two different replacements), and a global replace
# would not be able to distinguish them.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "731×6=4386"
$t.Cell(1,2).Range.Text  = "891×6=5346"
$t.Cell(1,3).Range.Text  = "691×6=4146"
$t.Cell(1,4).Range.Text  = "271×3=813"
$t.Cell(1,5).Range.Text  = "514×9=4626"

$t.Cell(5,1).Range.Text  = "134×9=1206"
$t.Cell(5,2).Range.Text  = "694×9=6246"
$t.Cell(5,3).Range.Text  = "496×5=2480"
$t.Cell(5,4).Range.Text  = "573×2=1146"
$t.Cell(5,5).Range.Text  = "150×6=900"

$t.Cell(10,1).Range.Text = "915×8=7320"
$t.Cell(10,2).Range.Text = "952×5=4760"
$t.Cell(10,3).Range.Text = "244×2=488"
$t.Cell(10,4).Range.Text = "182×5=910"
$t.Cell(10,5).Range.Text = "951×6=5706"

$t.Cell(15,1).Range.Text = "876×8=7008"
$t.Cell(15,2).Range.Text = "550×6=3300"
$t.Cell(15,3).Range.Text = "299×2=598"
$t.Cell(15,4).Range.Text = "902×6=5412"
$t.Cell(15,5).Range.Text = "882×7=6174"

$t.Cell(20,1).Range.Text = "752×5=3760"
$t.Cell(20,2).Range.Text = "844×6=5064"
$t.Cell(20,3).Range.Text = "976×5=4880"
$t.Cell(20,4).Range.Text = "410×6=2460"
$t.Cell(20,5).Range.Text = "888×3=2664"

Write-Output "done"
